# MAIN UPGRADE 2.7 to 3.8 - apply workbook edits
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("SOMINCOR_short")
$ws3 = $wb.Worksheets.Item("SOMINCOR_long")
$ws4 = $wb.Worksheets.Item("SOMINCOR_long_1")

# --- Sheet1: add 4 new data rows (22-25) ---
$newRows = @(
    @{ Row = 22; A = "SOMINCOR_MM_HEADS32_long";   B = 3696;  C = 270; D = 260; E = 4;  F = 1; G = 32; Fill = $true },
    @{ Row = 23; A = "SOMINCOR_MM_HEADS32_d_long"; B = 13171; C = 270; D = 260; E = 4;  F = 1; G = 32; Fill = $true },
    @{ Row = 24; A = "SOMINCOR_h_MF_long";         B = 3696;  C = 270; D = 260; E = 10; F = 1; G = 32; Fill = $false },
    @{ Row = 25; A = "SOMINCOR_h_MF_d_long";       B = 13171; C = 270; D = 260; E = 10; F = 1; G = 32; Fill = $false }
)

# Pre-seed the shared-string table in the same order the original author
# typed the labels (d_long variant before the long variant, per row pair)
# so that the regenerated sharedStrings.xml lines up with the source file.
$ws1.Range("A23").Value = "SOMINCOR_MM_HEADS32_d_long"
$ws1.Range("A22").Value = "SOMINCOR_MM_HEADS32_long"
$ws1.Range("A25").Value = "SOMINCOR_h_MF_d_long"
$ws1.Range("A24").Value = "SOMINCOR_h_MF_long"

foreach ($r in $newRows) {
    $row = $r.Row

    $cellA = $ws1.Range("A$row")
    $cellA.Borders.Color = 0
    $cellA.Borders.LineStyle = 1
    $cellA.Font.Bold = $true

    $cellB = $ws1.Range("B$row")
    $cellB.Borders.Color = 0
    $cellB.Borders.LineStyle = 1
    $cellB.Value = $r.B

    $cellC = $ws1.Range("C$row")
    $cellC.Borders.Color = 0
    $cellC.Borders.LineStyle = 1
    $cellC.Value = $r.C

    $cellD = $ws1.Range("D$row")
    $cellD.Borders.Color = 0
    $cellD.Borders.LineStyle = 1
    $cellD.Value = $r.D

    $cellE = $ws1.Range("E$row")
    $cellE.Borders.Color = 0
    $cellE.Borders.LineStyle = 1
    $cellE.Value = $r.E

    $cellF = $ws1.Range("F$row")
    $cellF.Borders.Color = 0
    $cellF.Borders.LineStyle = 1
    $cellF.Value = $r.F

    $cellG = $ws1.Range("G$row")
    $cellG.Borders.Color = 0
    $cellG.Borders.LineStyle = 1
    $cellG.Value = $r.G

    $cellH = $ws1.Range("H$row")
    $cellH.Borders.Color = 0
    $cellH.Borders.LineStyle = 1
    $cellH.Formula = "=B$row*C$row*D$row*E$row*F$row*G$row"

    $cellI = $ws1.Range("I$row")
    $cellI.Borders.Color = 0
    $cellI.Borders.LineStyle = 1
    $cellI.Formula = "=H$row/1000"

    $cellJ = $ws1.Range("J$row")
    $cellJ.Borders.Color = 0
    $cellJ.Borders.LineStyle = 1
    $cellJ.NumberFormat = "0.0"
    $cellJ.Formula = "=I$row/1000"

    $cellK = $ws1.Range("K$row")
    $cellK.Borders.Color = 0
    $cellK.Borders.LineStyle = 1
    $cellK.Font.Bold = $true
    $cellK.NumberFormat = "0.0"
    if ($r.Fill) {
        $cellK.Interior.Color = 255
    }
    $cellK.Formula = "=J$row/1000"
}

# --- SOMINCOR_short: clear the full-sheet selection ---
$ws2.Activate()
$ws2.Range("A1").Select()

# --- SOMINCOR_long: change selection ---
$ws3.Activate()
$ws3.Range("C4:C5").Select()

# --- SOMINCOR_long_1: clear selection / tabSelected ---
$ws4.Activate()
$ws4.Range("A1").Select()

# --- Sheet1: becomes the active tab / selection ---
$ws1.Activate()
$ws1.Range("B24").Select()
